# Refresh KHL referee stats (games played / PIM splits) and the as_of_utc
# timestamp for every referee row on the "Glavnye" (main) and "Lineynye" (linesmen) sheets.
$wb = $excel.ActiveWorkbook

# --- Worksheet #2: Главные ---
$ws = $wb.Worksheets.Item(2)

$rowUpdates = @(
    @{ Row = 2; Stats = @{ "C" = 26; "D" = 583; "E" = 239; "F" = 344; "G" = 22.42; "H" = 9.19; "I" = 13.23; "K" = 127; "L" = 5; "M" = 4 }; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 3; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 4; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 5; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 6; Stats = @{ "C" = 25; "D" = 437; "E" = 187; "F" = 250; "G" = 17.48; "H" = 7.48; "I" = 10; "J" = 86; "K" = 105 }; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 7; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 8; Stats = @{ "C" = 23; "D" = 416; "E" = 209; "F" = 207; "G" = 18.09; "H" = 9.09; "I" = 9; "J" = 97; "K" = 96 }; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 9; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 10; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 11; Stats = @{ "C" = 18; "D" = 470; "E" = 218; "F" = 252; "G" = 26.11; "H" = 12.11; "I" = 14; "J" = 94; "K" = 81; "L" = 4 }; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 12; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 13; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 14; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 15; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 16; Stats = @{ "C" = 25; "D" = 483; "E" = 238; "F" = 245; "G" = 19.32; "H" = 9.52; "I" = 9.800000000000001; "J" = 89; "K" = 90; "L" = 8 }; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 17; Stats = @{ "C" = 16; "D" = 264; "E" = 99; "F" = 165; "G" = 16.5; "H" = 6.19; "I" = 10.31; "J" = 47; "K" = 65; "L" = 1 }; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 18; Stats = @{ "C" = 24; "D" = 362; "E" = 170; "F" = 192; "G" = 15.08; "H" = 7.08; "I" = 8; "J" = 75; "K" = 91 }; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 19; Stats = @{ "C" = 20; "D" = 350; "E" = 168; "F" = 182; "G" = 17.5; "H" = 8.4; "I" = 9.1; "J" = 79; "K" = 76 }; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 20; Stats = @{ "C" = 24; "D" = 412; "E" = 174; "F" = 238; "G" = 17.17; "H" = 7.25; "I" = 9.92; "J" = 82; "K" = 89 }; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 21; Stats = @{ "C" = 21; "D" = 308; "E" = 138; "F" = 170; "G" = 14.67; "H" = 6.57; "J" = 59; "K" = 70 }; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 22; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 23; Stats = @{ "C" = 15; "D" = 196; "E" = 71; "F" = 125; "G" = 13.07; "H" = 4.73; "I" = 8.33; "J" = 33; "K" = 50; "L" = 1 }; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 24; Stats = @{ "C" = 26; "D" = 433; "E" = 193; "F" = 240; "G" = 16.65; "H" = 7.42; "I" = 9.23; "J" = 94; "K" = 105 }; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 25; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 26; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
)

foreach ($update in $rowUpdates) {
    $r = $update.Row
    foreach ($col in $update.Stats.Keys) {
        $ws.Range("$col$r").Value = $update.Stats[$col]
    }
    $ws.Range("AA$r").Value = $update.AsOfUtc
}

# --- Worksheet #3: Линейные ---
$ws = $wb.Worksheets.Item(3)

$rowUpdates = @(
    @{ Row = 2; Stats = @{ "C" = 16; "D" = 302; "E" = 131; "F" = 171; "G" = 18.88; "H" = 8.19; "I" = 10.69; "J" = 58; "K" = 63 }; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 3; Stats = @{ "C" = 24; "D" = 345; "E" = 171; "F" = 174; "G" = 14.38; "H" = 7.13; "I" = 7.25; "J" = 83; "K" = 72; "L" = 1 }; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 4; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 5; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 6; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 7; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 8; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 9; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 10; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 11; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 12; Stats = @{ "C" = 21; "D" = 372; "E" = 179; "F" = 193; "G" = 17.71; "H" = 8.52; "I" = 9.19; "J" = 82; "K" = 89 }; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 13; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 14; Stats = @{ "C" = 24; "D" = 402; "E" = 204; "F" = 198; "G" = 16.75; "H" = 8.5; "I" = 8.25; "J" = 102; "K" = 94 }; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 15; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 16; Stats = @{ "C" = 24; "D" = 420; "E" = 198; "F" = 222; "G" = 17.5; "H" = 8.25; "I" = 9.25; "K" = 96; "L" = 4; "M" = 6 }; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 17; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 18; Stats = @{ "C" = 27; "D" = 451; "E" = 214; "F" = 237; "G" = 16.7; "H" = 7.93; "I" = 8.779999999999999; "J" = 102; "K" = 101 }; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 19; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 20; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 21; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 22; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 23; Stats = @{ "C" = 14; "D" = 216; "E" = 103; "F" = 113; "G" = 15.43; "H" = 7.36; "I" = 8.07; "J" = 49; "K" = 54; "L" = 1 }; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 24; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 25; Stats = @{}; AsOfUtc = "2025-11-17 03:04:14" }
    @{ Row = 26; Stats = @{ "C" = 23; "D" = 461; "E" = 199; "F" = 262; "G" = 20.04; "H" = 8.65; "I" = 11.39; "J" = 77; "K" = 76 }; AsOfUtc = "2025-11-17 03:04:14" }
)

foreach ($update in $rowUpdates) {
    $r = $update.Row
    foreach ($col in $update.Stats.Keys) {
        $ws.Range("$col$r").Value = $update.Stats[$col]
    }
    $ws.Range("AA$r").Value = $update.AsOfUtc
}

